$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.35303951102532949
$ws.Range("A2").Value = -0.081372064400113686
$ws.Range("A3").Value = -0.0089999997028584744
$ws.Range("A4").Value = 0.061996874631255139
$ws.Range("A5").Value = -0.005999999708584447
$ws.Range("A6").Value = -0.0059999996995330207
$ws.Range("A7").Value = -0.019999999642971389
$ws.Range("A8").Value = -0.019999999638490529
$ws.Range("A9").Value = -0.0059999996881554551
$ws.Range("A10").Value = -0.0059999996834179115
$ws.Range("A11").Value = -0.0044999996888073213
$ws.Range("A12").Value = -0.005999999681663315
$ws.Range("A13").Value = -0.0059999996765451868
$ws.Range("A14").Value = -0.01199999965138332
$ws.Range("A15").Value = -0.0059999996738531181
$ws.Range("A16").Value = -0.00599999967279774
$ws.Range("A17").Value = -0.0059999996713617776
$ws.Range("A18").Value = -0.0089999996592649012
$ws.Range("A19").Value = -0.008999999708775519
$ws.Range("A20").Value = -0.0089999997020040468
$ws.Range("A21").Value = -0.0089999997006309229
$ws.Range("A22").Value = -0.044047235793157835
$ws.Range("A23").Value = -0.0089999997001806165
$ws.Range("A24").Value = -0.041999999565799584
$ws.Range("A25").Value = -0.041999999563474333
$ws.Range("A26").Value = -0.0059999996995010463
$ws.Range("A27").Value = -0.0059999996990649507
$ws.Range("A28").Value = -0.0059999996976367598
$ws.Range("A29").Value = -0.011999999673374617
$ws.Range("A30").Value = -0.019999999641928223
$ws.Range("A31").Value = -0.014999999662880015
$ws.Range("A32").Value = 0.016721326461294872
$ws.Range("A33").Value = 0.017392593645193521
